# Personal_Fit_Model/nyc_race.xlsx edit
#
# 1. Remove the "BX98 Rikers Island" row (row 57) entirely - no work force there.
# 2. Simplify the "MN01 Marble Hill2-Inwood" rich-text label (with superscript "2")
#    down to plain text "MN01 Marble Hill-Inwood" (this row shifts up to row 115
#    once the Rikers Island row above it is removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select then delete the whole "BX98 Rikers Island" row (row 57). Selecting first
# mirrors how this was actually done in the UI and leaves the selection sitting on
# the row that slides up into its place.
$ws.Rows("57").Select()
$ws.Rows("57").Delete()

# After the deletion above, the "MN01 Marble Hill2-Inwood" entry (originally row
# 116) has shifted up to row 115. Replace its superscripted rich text with plain
# text that drops the footnote number.
$ws.Range("A115").Value = "MN01 Marble Hill-Inwood"
